# Append the new portfolio row (2025-09-08) as row 24, extending the
# data range from A1:D23 to A1:D24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric columns first (these are plain numbers, no special formatting).
$ws.Range("B24").Value = 57.45000076293945
$ws.Range("C24").Value = 719.5
$ws.Range("D24").Value = 329.6499938964844

# Date column: stored as literal text (matching the rest of column A),
# not as a date serial number. Temporarily force a text number format so
# Excel doesn't auto-convert the "YYYY-MM-DD" string into a date value,
# then restore the cell to the default/unstyled state.
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = "2025-09-08"
$ws.Range("A24").Style = "Normal"
